# Applies the weekly cryptos price/volume refresh described in the commit
# "Updated cryptos list on Wed Jun  5 22:36:09 UTC 2024 with GitHub Actions".
# Columns D (Price) and E (Volume(1h)) are plain text cells (coinranking.com
# renders thousands-separated prices and padded +/-% strings), so every write
# goes through .Value as a string. A handful of Price cells (e.g. "696.12")
# parse as plain numbers; for those we briefly force a Text number format so
# Excel keeps storing/display them as text (matching the source feed), then
# restore "General" so no lasting number-format change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.171.42'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '3.847.21'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '696.12'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.68'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +2.18%  '
$ws.Range('D7').Value = '3.843.46'
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E13').Value = '  +4.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.44'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('D15').Value = '4.492.58'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').Value = '3.854.67'
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('D17').Value = '71.195.52'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.71'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.19'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '493.10'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +3.26%  '
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.00'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000145'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +2.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.32'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.60'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +2.68%  '
$ws.Range('E28').Value = '  +2.01%  '
$ws.Range('D29').Value = '4.004.68'
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.18'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +7.44%  '
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.68'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.28'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +1.57%  '
$ws.Range('D37').Value = '3.796.59'
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +2.43%  '
$ws.Range('E40').Value = '  +12.43%  '
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.03'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +6.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '163.32'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +2.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000307'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +3.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '48.64'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '44.21'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -4.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '418.93'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +5.63%  '
$ws.Range('E51').Value = '  +1.06%  '
